$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record needs to be inserted at row 32 (just below the
# header + the first 30 data rows), pushing the existing rows 32..113 down
# to 33..114 and extending the used range from A1:R113 to A1:R114.
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with the new record's data.
$ws.Range("A32").Value = 11
$ws.Range("B32").Value = "Vega Monumental Concepción"
$ws.Range("C32").Value = "Bíobío"
$ws.Range("D32").Value = 44804
$ws.Range("E32").Value = 8
$ws.Range("F32").Value = 100112001
$ws.Range("G32").Value = "Berenjena"
$ws.Range("H32").Value = "Sin especificar"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 130
$ws.Range("K32").Value = 11000
$ws.Range("L32").Value = 12000
$ws.Range("M32").Value = 11538
$ws.Range("N32").Value = "$/caja 60 unidades"
$ws.Range("O32").Value = "Región de Arica y Parinacota"
$ws.Range("P32").Value = 192
$ws.Range("Q32").Value = 60
$ws.Range("R32").Value = "Hortaliza"
